$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: Expand the description of the fuzzy-match step (the paragraph
# that begins "If the disease name contained one of the above tumor key
# words..."). The old text trailed off with "...of 10% from  "; the new
# text spells out the exact-match check plus the fuzzy (20%) check and
# ends with "...it was flagged as a potential tumor."
# ---------------------------------------------------------------------
$oldPara1 = "If the disease name contained one of the above tumor key words, we flagged that disease as a potential tumor.   In our second step in detecting tumors, we considered the tumor names contained in the WHO database (5 edition) and developed a fuzzy string match program to match disease names to tumor names in the WHO database. If a disease from clinical trials was within an edit distance (edits such as deletion, insertion or substitutions needed to convert one string to another string) of 10% from  "

$newPara1 = "If the disease name contained one of the above tumor key words, we flagged that disease as a potential tumor.  In our second step in detecting tumors, we considered the tumor names contained in the WHO database (5 edition) and developed a fuzzy string match program to match disease names to tumor names in the WHO database. If a disease from clinical trials exactly matched a term in the WHO database, it was flagged as tumor. If the disease did not match to any tumor within the WHO database, we performed a fuzzy (approximate) match with the disease name with each term in the WHO database. This was done by computing the edit distance (edits such as deletion, insertion or substitutions needed to convert one string to another string) of the clinical trial disease to each WHO database term and then if there were any matches within 20% any WHO tumor names then it was flagged as a potential tumor."

$d.Content.Find.Execute($oldPara1, $true, $false, $false, $false, $false, $true, 1, $false, $newPara1, 2) | Out-Null

# ---------------------------------------------------------------------
# Step 2: The following paragraph used to start with the now-redundant
# "any WHO tumor names then it was flagged as a potential tumor." lead-in
# (that sentence now lives at the end of paragraph 1, above). Locate that
# paragraph, strip the stale lead-in sentence from its front, and push
# the remaining "  After these two steps..." text out into its own new
# paragraph - matching how the sentence was split off in the real edit.
# ---------------------------------------------------------------------
$leadIn = "any WHO tumor names then it was flagged as a potential tumor."

$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "$leadIn*") {
        $targetIndex = $i
        break
    }
}

$paraRange = $d.Paragraphs.Item($targetIndex).Range
$paraRange.Find.Execute($leadIn, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# $paraRange is now collapsed to the matched "lead-in" text (it sits right
# at the start of the paragraph). Clear it, then drop in a paragraph
# break so the remaining text becomes a paragraph of its own.
$paraRange.Text = ""
$paraRange.InsertParagraphBefore() | Out-Null
